$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.767.07'
$ws.Range('E2').Value = '  +0.47%  '
$ws.Range('D3').Value = '2.672.78'
$ws.Range('E3').Value = '  +0.81%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '602.08'
$ws.Range('E5').Value = '  -0.86%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.94'
$ws.Range('E6').Value = '  +0.57%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.617'
$ws.Range('E8').Value = '  +5.03%  '
$ws.Range('E9').Value = '  +1.31%  '
$ws.Range('E10').Value = '  +1.13%  '
$ws.Range('E11').Value = '  -0.39%  '
$ws.Range('E12').Value = '  -0.27%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '29.46'
$ws.Range('E13').Value = '  -1.53%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000197'
$ws.Range('E14').Value = '  +0.78%  '
$ws.Range('D15').Value = '3.154.49'
$ws.Range('E15').Value = '  +0.88%  '
$ws.Range('D16').Value = '65.573.84'
$ws.Range('E16').Value = '  +0.53%  '
$ws.Range('D17').Value = '2.666.53'
$ws.Range('E17').Value = '  +0.39%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.63'
$ws.Range('E18').Value = '  -0.80%  '
$ws.Range('E19').Value = '  -1.27%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.58'
$ws.Range('E20').Value = '  +1.74%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '352.22'
$ws.Range('E21').Value = '  -1.69%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.21%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.88'
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000110'
$ws.Range('E24').Value = '  +4.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.79'
$ws.Range('E25').Value = '  +3.91%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.62'
$ws.Range('E26').Value = '  -4.73%  '
$ws.Range('E27').Value = '  +1.11%  '
$ws.Range('E28').Value = '  -1.61%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.13'
$ws.Range('E29').Value = '  +0.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '543.53'
$ws.Range('E30').Value = '  +3.93%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.14'
$ws.Range('E32').Value = '  -1.13%  '
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('E34').Value = '  +3.70%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.44'
$ws.Range('E35').Value = '  -1.18%  '
$ws.Range('E36').Value = '  -1.95%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.41'
$ws.Range('E37').Value = '  -1.22%  '
$ws.Range('E38').Value = '  -0.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '158.98'
$ws.Range('E39').Value = '  -2.30%  '
$ws.Range('E40').Value = '  -1.13%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('E42').Value = '  +2.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '164.87'
$ws.Range('E43').Value = '  -0.36%  '
$ws.Range('E44').Value = '  -0.68%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0613'
$ws.Range('E45').Value = '  +0.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.30'
$ws.Range('E46').Value = '  -1.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '23.22'
$ws.Range('E47').Value = '  +1.18%  '
$ws.Range('E48').Value = '  -0.89%  '
$ws.Range('E49').Value = '  -1.39%  '
$ws.Range('E50').Value = '  +2.88%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '20.27'
$ws.Range('E51').Value = '  +3.06%  '
